$wb = $excel.ActiveWorkbook

# Duplicate the last tire sheet (Tir_145_70R13) to create a new sheet for the
# 430/50R38 tire, placing it after the existing sheets.
$src = $wb.Worksheets.Item("Tir_145_70R13")
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Tir_430_50R38"

# Update the identifying label and the key tire dimensions/mass for this new
# tire size. Convention change to support multi-axle vehicles.
$ws.Range("H3").Value = "Testrig_Post_430_50R38"
$ws.Range("H5").Value = 0.6731
$ws.Range("H6").Value = 0.47752
$ws.Range("H7").Value = 0.42951
$ws.Range("H7").NumberFormat = "0.000"
$ws.Range("H8").Value = 0.41

$ws.Activate()
